# Refactoring out query strings.
# Replace the placeholder test data ("TESET"/"TEST") in row 5 of the
# JudgeDashboard sheet with real defendant data ("DOUGLAS"/"MICHAEL"),
# and move the active selection to the cell that was edited (D5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C = Lastname, Column D = Firstname
$ws.Range("C5").Value = "DOUGLAS"
$ws.Range("D5").Value = "MICHAEL"

# Update the active selection to reflect the edited cell.
$ws.Range("D5").Select()
